$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.938.95'
$ws.Range("E2").Value = '  +0.15%  '
$ws.Range("D3").Value = '2.214.45'
$ws.Range("E3").Value = '  -1.82%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '230.80'
$ws.Range("E5").Value = '  -0.14%  '
$ws.Range("E6").Value = '  -1.78%  '
$ws.Range("D7").Value = '60.72'
$ws.Range("E7").Value = '  -2.12%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = '0.402'
$ws.Range("E9").Value = '  -0.96%  '
$ws.Range("D10").Value = '57.20'
$ws.Range("D11").Value = '0.0899'
$ws.Range("E11").Value = '  +1.47%  '
$ws.Range("E12").Value = '  -1.03%  '
$ws.Range("D13").Value = '2.543.52'
$ws.Range("E13").Value = '  -1.73%  '
$ws.Range("D14").Value = '15.43'
$ws.Range("E14").Value = '  -3.85%  '
$ws.Range("D15").Value = '22.09'
$ws.Range("E15").Value = '  +0.69%  '
$ws.Range("D16").Value = '0.797'
$ws.Range("E16").Value = '  -1.68%  '
$ws.Range("E17").Value = '  -1.67%  '
$ws.Range("D18").Value = '2.203.71'
$ws.Range("E18").Value = '  -2.14%  '
$ws.Range("D19").Value = '41.853.90'
$ws.Range("E19").Value = '  +0.26%  '
$ws.Range("D20").Value = '0.0₃0923'
$ws.Range("E20").Value = '  +2.19%  '
$ws.Range("D21").Value = '72.02'
$ws.Range("E21").Value = '  -2.22%  '
$ws.Range("D22").Value = '6.05'
$ws.Range("E22").Value = '  -2.04%  '
$ws.Range("D23").Value = '242.59'
$ws.Range("E23").Value = '  -2.43%  '
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("D25").Value = '2.38'
$ws.Range("E25").Value = '  +1.79%  '
$ws.Range("D26").Value = '2.34'
$ws.Range("E26").Value = '  -2.34%  '
$ws.Range("D27").Value = '9.61'
$ws.Range("E27").Value = '  -1.64%  '
$ws.Range("D28").Value = '168.99'
$ws.Range("E28").Value = '  +0.12%  '
$ws.Range("E29").Value = '  -3.35%  '
$ws.Range("D30").Value = '1.43'
$ws.Range("E30").Value = '  -2.56%  '
$ws.Range("D31").Value = '19.78'
$ws.Range("E31").Value = '  -2.48%  '
$ws.Range("E32").Value = '  -5.85%  '
$ws.Range("E34").Value = '  -4.22%  '
$ws.Range("D35").Value = '4.60'
$ws.Range("E35").Value = '  -2.15%  '
$ws.Range("D36").Value = '0.0650'
$ws.Range("E36").Value = '  +2.85%  '
$ws.Range("E37").Value = '  -6.51%  '
$ws.Range("E38").Value = '  -8.16%  '
$ws.Range("E39").Value = '  -4.22%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '0.0243'
$ws.Range("E40").Value = '  +1.39%  '
$ws.Range("B41").Value = 'BinanceUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.42%  '
$ws.Range("B42").Value = 'TerraClassic'
$ws.Range("C42").Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range("D42").Value = '0.000233'
$ws.Range("E42").Value = '  -5.79%  '
$ws.Range("D43").Value = '8.48'
$ws.Range("E43").Value = '  -3.27%  '
$ws.Range("D44").Value = '0.0952'
$ws.Range("E44").Value = '  -2.38%  '
$ws.Range("E45").Value = '  -0.79%  '
$ws.Range("D46").Value = '96.73'
$ws.Range("E46").Value = '  -4.44%  '
$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").Value = '4.31'
$ws.Range("E47").Value = '  -11.98%  '
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '1.450.51'
$ws.Range("E48").Value = '  -2.63%  '
$ws.Range("E49").Value = '  -1.61%  '
$ws.Range("D50").Value = '16.10'
$ws.Range("E50").Value = '  -4.01%  '
$ws.Range("D51").Value = '1.06'
$ws.Range("E51").Value = '  -4.09%  '
